$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 805 (the "ハム太郎って子を知ってるかい" post) — all
# subsequent rows shift up by one to close the gap.
$ws.Rows("805").Delete()
